$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet is a daily log of price observations for "Pepino ensalada" at
# "Feria Lagunitas de Puerto Montt". A new (most recent) observation is
# being recorded, which pushes every existing row from 261 downward down
# by one (last row 341 becomes 342).
$ws.Rows.Item(261).Insert()

$ws.Range("A261").Value = 4
$ws.Range("B261").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C261").Value = "Los Lagos"
$ws.Range("D261").Value = 44876
$ws.Range("E261").Value = 10
$ws.Range("F261").Value = 100112043
$ws.Range("G261").Value = "Pepino ensalada"
$ws.Range("H261").Value = "Sin especificar"
$ws.Range("I261").Value = "Primera"
$ws.Range("J261").Value = 400
$ws.Range("K261").Value = 25000
$ws.Range("L261").Value = 25000
$ws.Range("M261").Value = 25000
$ws.Range("N261").Value = "`$/caja 60 unidades"
$ws.Range("O261").Value = "Región de Arica y Parinacota"
$ws.Range("P261").Value = 417
$ws.Range("Q261").Value = 60
$ws.Range("R261").Value = "Hortaliza"
